$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.474131
$ws.Range("H2").Value = 14.948262
$ws.Range("I2").Value = 0.002368493638035435
$ws.Range("J2").Value = 0.001580807736380949
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3657716666666667
$ws.Range("N2").Value = 1.097315
$ws.Range("O2").Value = 0.03059585711603819
$ws.Range("P2").Value = 0.03059585711603819
$ws.Range("Q2").Value = 2.733825352755
$ws.Range("R2").Value = 16.40295211653
$ws.Range("S2").Value = 0.00007246609292957766
$ws.Range("T2").Value = 0.00004836616763023926

$ws.Range("G3").Value = 7.474131
$ws.Range("H3").Value = 14.948262
$ws.Range("I3").Value = 0.002368493638035435
$ws.Range("J3").Value = 0.001580807736380949
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.510206
$ws.Range("N3").Value = 34.530618
$ws.Range("O3").Value = 0.9627990635838353
$ws.Range("P3").Value = 0.9627990635838353
$ws.Range("Q3").Value = 86.02878748098598
$ws.Range("R3").Value = 516.172724885916
$ws.Range("S3").Value = 0.002280383456804788
$ws.Range("T3").Value = 0.00152200020829366

$ws.Range("G4").Value = 7.474131
$ws.Range("H4").Value = 14.948262
$ws.Range("I4").Value = 0.002368493638035435
$ws.Range("J4").Value = 0.001580807736380949
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.07896333333333333
$ws.Range("N4").Value = 0.23689
$ws.Range("O4").Value = 0.006605079300126477
$ws.Range("P4").Value = 0.006605079300126477
$ws.Range("Q4").Value = 0.59018229753
$ws.Range("R4").Value = 3.54109378518
$ws.Range("S4").Value = 0.00001564408830106911
$ws.Range("T4").Value = 0.0000104413604570496

$ws.Range("G5").Value = 28.63660866666666
$ws.Range("H5").Value = 85.909826
$ws.Range("I5").Value = 0.009074717240293273
$ws.Range("J5").Value = 0.009085130938428906
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3657716666666667
$ws.Range("N5").Value = 1.097315
$ws.Range("O5").Value = 0.03059585711603819
$ws.Range("P5").Value = 0.03059585711603819
$ws.Range("Q5").Value = 10.47446007968778
$ws.Range("R5").Value = 94.27014071719
$ws.Range("S5").Value = 0.0002776487520524613
$ws.Range("T5").Value = 0.0002779673680726688

$ws.Range("G6").Value = 28.63660866666666
$ws.Range("H6").Value = 85.909826
$ws.Range("I6").Value = 0.009074717240293273
$ws.Range("J6").Value = 0.009085130938428906
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.510206
$ws.Range("N6").Value = 34.530618
$ws.Range("O6").Value = 0.9627990635838353
$ws.Range("P6").Value = 0.9627990635838353
$ws.Range("Q6").Value = 329.6132648947186
$ws.Range("R6").Value = 2966.519384052468
$ws.Range("S6").Value = 0.00873712926124245
$ws.Range("T6").Value = 0.008747155560055881

$ws.Range("G7").Value = 28.63660866666666
$ws.Range("H7").Value = 85.909826
$ws.Range("I7").Value = 0.009074717240293273
$ws.Range("J7").Value = 0.009085130938428906
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.07896333333333333
$ws.Range("N7").Value = 0.23689
$ws.Range("O7").Value = 0.006605079300126477
$ws.Range("P7").Value = 0.006605079300126477
$ws.Range("Q7").Value = 2.261242075682222
$ws.Range("R7").Value = 20.35117868114
$ws.Range("S7").Value = 0.00005993922699836196
$ws.Range("T7").Value = 0.0000600080103003554

$ws.Range("G8").Value = 403.4856263333333
$ws.Range("H8").Value = 1210.456879
$ws.Range("I8").Value = 0.1278614382072301
$ws.Range("J8").Value = 0.1280081656903483
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3657716666666667
$ws.Range("N8").Value = 1.097315
$ws.Range("O8").Value = 0.03059585711603819
$ws.Range("P8").Value = 0.03059585711603819
$ws.Range("Q8").Value = 147.5836100199872
$ws.Range("R8").Value = 1328.252490179885
$ws.Range("S8").Value = 0.003912030294039558
$ws.Range("T8").Value = 0.00391651954714804

$ws.Range("G9").Value = 403.4856263333333
$ws.Range("H9").Value = 1210.456879
$ws.Range("I9").Value = 0.1278614382072301
$ws.Range("J9").Value = 0.1280081656903483
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.510206
$ws.Range("N9").Value = 34.530618
$ws.Range("O9").Value = 0.9627990635838353
$ws.Range("P9").Value = 0.9627990635838353
$ws.Range("Q9").Value = 4644.20267713569
$ws.Range("R9").Value = 41797.82409422121
$ws.Range("S9").Value = 0.1231048729744035
$ws.Range("T9").Value = 0.1232461420577518

$ws.Range("G10").Value = 403.4856263333333
$ws.Range("H10").Value = 1210.456879
$ws.Range("I10").Value = 0.1278614382072301
$ws.Range("J10").Value = 0.1280081656903483
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.07896333333333333
$ws.Range("N10").Value = 0.23689
$ws.Range("O10").Value = 0.006605079300126477
$ws.Range("P10").Value = 0.006605079300126477
$ws.Range("Q10").Value = 31.86057000736777
$ws.Range("R10").Value = 286.74513006631
$ws.Range("S10").Value = 0.0008445349387869761
$ws.Range("T10").Value = 0.0008455040854484802

$ws.Range("G11").Value = 3.377213
$ws.Range("H11").Value = 6.754426
$ws.Range("I11").Value = 0.001070212377170077
$ws.Range("J11").Value = 0.0007142936667562174
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3657716666666667
$ws.Range("N11").Value = 1.097315
$ws.Range("O11").Value = 0.03059585711603819
$ws.Range("P11").Value = 0.03059585711603819
$ws.Range("Q11").Value = 1.235288827698333
$ws.Range("R11").Value = 7.411732966190001
$ws.Range("S11").Value = 0.00003274406497571126
$ws.Range("T11").Value = 0.00002185442696696422

$ws.Range("G12").Value = 3.377213
$ws.Range("H12").Value = 6.754426
$ws.Range("I12").Value = 0.001070212377170077
$ws.Range("J12").Value = 0.0007142936667562174
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 11.510206
$ws.Range("N12").Value = 34.530618
$ws.Range("O12").Value = 0.9627990635838353
$ws.Range("P12").Value = 0.9627990635838353
$ws.Range("Q12").Value = 38.872417335878
$ws.Range("R12").Value = 233.234504015268
$ws.Range("S12").Value = 0.001030399474575181
$ws.Range("T12").Value = 0.0006877212734767501

$ws.Range("G13").Value = 3.377213
$ws.Range("H13").Value = 6.754426
$ws.Range("I13").Value = 0.001070212377170077
$ws.Range("J13").Value = 0.0007142936667562174
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.07896333333333333
$ws.Range("N13").Value = 0.23689
$ws.Range("O13").Value = 0.006605079300126477
$ws.Range("P13").Value = 0.006605079300126477
$ws.Range("Q13").Value = 0.2666759958566667
$ws.Range("R13").Value = 1.60005597514
$ws.Range("S13").Value = 0.000007068837619185228
$ws.Range("T13").Value = 0.000004717966312502931

$ws.Range("G14").Value = 2406.316202666666
$ws.Range("H14").Value = 7218.948608
$ws.Range("I14").Value = 0.7625427781661288
$ws.Range("J14").Value = 0.7634178346661893
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.3657716666666667
$ws.Range("N14").Value = 1.097315
$ws.Range("O14").Value = 0.03059585711603819
$ws.Range("P14").Value = 0.03059585711603819
$ws.Range("Q14").Value = 880.162287976391
$ws.Range("R14").Value = 7921.46059178752
$ws.Range("S14").Value = 0.02333064988563768
$ws.Range("T14").Value = 0.02335742298928199

$ws.Range("G15").Value = 2406.316202666666
$ws.Range("H15").Value = 7218.948608
$ws.Range("I15").Value = 0.7625427781661288
$ws.Range("J15").Value = 0.7634178346661893
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 11.510206
$ws.Range("N15").Value = 34.530618
$ws.Range("O15").Value = 0.9627990635838353
$ws.Range("P15").Value = 0.9627990635838353
$ws.Range("Q15").Value = 27697.19519383108
$ws.Range("R15").Value = 249274.7567444797
$ws.Range("S15").Value = 0.7341754727609651
$ws.Range("T15").Value = 0.7350179763398063

$ws.Range("G16").Value = 2406.316202666666
$ws.Range("H16").Value = 7218.948608
$ws.Range("I16").Value = 0.7625427781661288
$ws.Range("J16").Value = 0.7634178346661893
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.07896333333333333
$ws.Range("N16").Value = 0.23689
$ws.Range("O16").Value = 0.006605079300126477
$ws.Range("P16").Value = 0.006605079300126477
$ws.Range("Q16").Value = 190.0107484165689
$ws.Range("R16").Value = 1710.09673574912
$ws.Range("S16").Value = 0.005036655519526033
$ws.Range("T16").Value = 0.005042435337101025

$ws.Range("G17").Value = 306.3577066666667
$ws.Range("H17").Value = 919.07312
$ws.Range("I17").Value = 0.09708236037114229
$ws.Range("J17").Value = 0.0971937673018961
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3657716666666667
$ws.Range("N17").Value = 1.097315
$ws.Range("O17").Value = 0.03059585711603819
$ws.Range("P17").Value = 0.03059585711603819
$ws.Range("Q17").Value = 112.0569689636445
$ws.Range("R17").Value = 1008.5127206728
$ws.Range("S17").Value = 0.002970318026403197
$ws.Range("T17").Value = 0.002973726616938278

$ws.Range("G18").Value = 306.3577066666667
$ws.Range("H18").Value = 919.07312
$ws.Range("I18").Value = 0.09708236037114229
$ws.Range("J18").Value = 0.0971937673018961
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 11.510206
$ws.Range("N18").Value = 34.530618
$ws.Range("O18").Value = 0.9627990635838353
$ws.Range("P18").Value = 0.9627990635838353
$ws.Range("Q18").Value = 3526.240313420906
$ws.Range("R18").Value = 31736.16282078816
$ws.Range("S18").Value = 0.09347080565584423
$ws.Range("T18").Value = 0.09357806814445076

$ws.Range("G19").Value = 306.3577066666667
$ws.Range("H19").Value = 919.07312
$ws.Range("I19").Value = 0.09708236037114229
$ws.Range("J19").Value = 0.0971937673018961
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.07896333333333333
$ws.Range("N19").Value = 0.23689
$ws.Range("O19").Value = 0.006605079300126477
$ws.Range("P19").Value = 0.006605079300126477
$ws.Range("Q19").Value = 24.19102571075555
$ws.Range("R19").Value = 217.7192313968
$ws.Range("S19").Value = 0.000641236688894851
$ws.Range("T19").Value = 0.0006419725405070636

